$wb = $excel.ActiveWorkbook

# --- 1. Status text: "Ready for handoff" -> "In Translation" --------------
# This text lives in the shared-string table and is referenced by every
# cell that shows it: Overview!E2, Overview!F2, zh-cn!C2 and de-de!C2.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Column widths: 17.2159881591797 -> 13.4101845877511 ---------------
# Overview columns E and F, and column C ("Status") on the zh-cn / de-de
# sheets, get narrower. ColumnWidth closest achievable setting.
$wsOverview.Columns("E").ColumnWidth = 12.5
$wsOverview.Columns("F").ColumnWidth = 12.5

$wsZhCn.Columns("C").ColumnWidth = 12.5

$wsDeDe.Columns("C").ColumnWidth = 12.5
